$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet -- the whole table is being restructured
# (two header rows collapse into one, new idx/idx2/Name/Date columns are
# added in front, and the data block shifts up by one row).
$ws.Cells.Clear()

# ---------------------------------------------------------------------
# Power-plant master data (name, canton id, years in service, capacities)
# ---------------------------------------------------------------------
$names      = @("Rotzloch","Obermatt","Oberrickenbach","Wolfenschiessen","Sustli","Dallenwil","Obermatt-Nebenzentrale","Arni","Engelberg")
$idx2       = @(304100,303000,303300,303400,302600,303200,303100,302900,302800)
$dateStart  = @(1872,1905,1937,1945,1957,1962,1963,1966,1967)
$dateEnd    = @(1935,1963,1991,1983,1998,1987,$null,$null,$null)
$m3s        = @(0.6,11,1,2.6,0.51,14.7,11,1,1.4)
$mw1        = @(0.3,7.73,8.7,6.6,1.63,12.22,0.48,1.85,7.74)
$mw2        = @(0.3,7.19,6.8,6.6,1.56,11.06,0.44,1.74,7.64)
$gwhWinter  = @(0.75,6.29,4.4,6.6,3.5,13.13,0.24,0.29,2.91)
$gwhSummer  = @(0.75,23.96,9.8,13.7,5.3,40.67,0.97,4.06,14.65)
$gwhYear    = @(1.5,30.25,14.2,20.3,8.8,53.8,1.21,4.35,17.56)

# ---------------------------------------------------------------------
# Write cells in the exact order the strings are first used so the
# shared-string table comes out in the right sequence: (m3/s) first,
# then the plant names, then the new column headers.
# ---------------------------------------------------------------------

# (m3/s) unit cell first -> shared string index 0
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("F1").Font.Size = 9

# Plant names (rows 2..10) -> shared string indices 1..9
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $names[$i]
    $ws.Cells.Item($row, 3).Font.Size = 9
}

# New header labels for columns A-E
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# New header labels for columns G-K
$ws.Range("G1").Value = "(MW1)"
$ws.Range("G1").Font.Size = 9
$ws.Range("H1").Value = "(MW2)"
$ws.Range("H1").Font.Size = 9
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("I1").Font.Size = 9
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("J1").Font.Size = 9
$ws.Range("K1").Value = "(GWh) Year"
$ws.Range("K1").Font.Size = 9

# ---------------------------------------------------------------------
# Data rows 2..10
# ---------------------------------------------------------------------
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 1).Font.Size = 9
    $ws.Cells.Item($row, 1).NumberFormat = "0"

    $ws.Cells.Item($row, 2).Value = $idx2[$i]
    $ws.Cells.Item($row, 2).Font.Size = 9
    $ws.Cells.Item($row, 2).NumberFormat = "0"

    $ws.Cells.Item($row, 4).Value = $dateStart[$i]
    $ws.Cells.Item($row, 4).Font.Size = 9
    $ws.Cells.Item($row, 4).NumberFormat = "0"

    if ($dateEnd[$i] -ne $null) {
        $ws.Cells.Item($row, 5).Value = $dateEnd[$i]
        $ws.Cells.Item($row, 5).Font.Size = 9
        $ws.Cells.Item($row, 5).NumberFormat = "0"
    }

    $ws.Cells.Item($row, 6).Value = $m3s[$i]
    $ws.Cells.Item($row, 6).Font.Size = 9
    $ws.Cells.Item($row, 6).NumberFormat = "0.00"

    $ws.Cells.Item($row, 7).Value = $mw1[$i]
    $ws.Cells.Item($row, 7).Font.Size = 9
    $ws.Cells.Item($row, 7).NumberFormat = "0.00"

    $ws.Cells.Item($row, 8).Value = $mw2[$i]
    $ws.Cells.Item($row, 8).Font.Size = 9
    $ws.Cells.Item($row, 8).NumberFormat = "0.00"

    $ws.Cells.Item($row, 9).Value = $gwhWinter[$i]
    $ws.Cells.Item($row, 9).Font.Size = 9
    $ws.Cells.Item($row, 9).NumberFormat = "0.00"

    $ws.Cells.Item($row, 10).Value = $gwhSummer[$i]
    $ws.Cells.Item($row, 10).Font.Size = 9
    $ws.Cells.Item($row, 10).NumberFormat = "0.00"

    $ws.Cells.Item($row, 11).Value = $gwhYear[$i]
    $ws.Cells.Item($row, 11).Font.Size = 9
    $ws.Cells.Item($row, 11).NumberFormat = "0.00"
}

# Match the workbook's recorded selection after the edit.
$ws.Range("A2:K2").Select()
